$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("M2:M4").NumberFormat = "@"

$ws.Range("M2").Value = "2022-01-12"
$ws.Range("M3").Value = "2022-01-13"
$ws.Range("M4").Value = "2022-01-14"

$ws.Range("M4").Select()
